$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("error message")

# Insert two new rows before row 24, shifting the old rows 24-25 down to 26-27
$ws.Rows.Item(24).Resize(2).Insert()

# Fill in the new rows 24 and 25
$ws.Range("A24").Value = "E   "
$ws.Range("C24").Value = "beas_qc_1006"
$ws.Range("E24").Value = "S/N entered cannot be duplicated in one QC order"

$ws.Range("A25").Value = "E   "
$ws.Range("C25").Value = "beas_qc_1007"
$ws.Range("E25").Value = "S/N is mandatory."

# Update the selection to match the target state
$ws.Range("E13").Select()
